$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "File" column (E) on row 3 references an existing derivative fixture
# instead of a made-up name.
$ws.Range("E3").Value = "videoshort.mp4"

# Reflect the cell the author was on when making the edit.
$ws.Range("E3").Select()
